$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44307
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'

# Row 3
$ws.Range("D3").Value = 44307
$ws.Range("Q3").Value = '$/bandeja 18 kilos granel'

# Row 4
$ws.Range("D4").Value = 44285
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 9000
$ws.Range("P4").Value = 9500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("S4").Value = 528
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44285
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("S5").Value = 444
$ws.Range("T5").Value = 18

# Row 8
$ws.Range("D8").Value = 44202
$ws.Range("K8").Value = 'Black Amber'
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("S8").Value = 806
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44202
$ws.Range("K9").Value = 'Black Amber'
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44323
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 11000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 11500
$ws.Range("S10").Value = 639

# Row 11
$ws.Range("D11").Value = 44323
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("S11").Value = 500

# Row 12
$ws.Range("D12").Value = 44189
$ws.Range("K12").Value = 'Red Beaut'
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 12500
$ws.Range("Q12").Value = '$/caja 15 kilos granel'
$ws.Range("S12").Value = 833
$ws.Range("T12").Value = 15

# Row 13
$ws.Range("D13").Value = 44189
$ws.Range("K13").Value = 'Red Beaut'
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("Q13").Value = '$/caja 15 kilos granel'
$ws.Range("S13").Value = 667
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 44246
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("Q14").Value = '$/caja 16 kilos granel'
$ws.Range("S14").Value = 625
$ws.Range("T14").Value = 16

# Row 15
$ws.Range("D15").Value = 44246
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("Q15").Value = '$/caja 16 kilos granel'
$ws.Range("T15").Value = 16

# Row 18
$ws.Range("D18").Value = 44328
$ws.Range("K18").Value = 'Angeleno'
$ws.Range("M18").Value = 100

# Row 19
$ws.Range("D19").Value = 44328
$ws.Range("K19").Value = 'Angeleno'
$ws.Range("M19").Value = 50

# Row 20
$ws.Range("D20").Value = 44343
$ws.Range("K20").Value = 'Angeleno'
$ws.Range("Q20").Value = '$/bandeja 18 kilos granel'
$ws.Range("S20").Value = 583
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("D21").Value = 44343
$ws.Range("K21").Value = 'Angeleno'
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 9000
$ws.Range("P21").Value = 9000
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
$ws.Range("T21").Value = 18

# Row 24
$ws.Range("D24").Value = 44218
$ws.Range("K24").Value = 'Black Amber'
$ws.Range("N24").Value = 10000
$ws.Range("O24").Value = 11000
$ws.Range("P24").Value = 10500
$ws.Range("Q24").Value = '$/caja 16 kilos granel'
$ws.Range("S24").Value = 656
$ws.Range("T24").Value = 16

# Row 25
$ws.Range("D25").Value = 44218
$ws.Range("K25").Value = 'Black Amber'
$ws.Range("Q25").Value = '$/caja 16 kilos granel'
$ws.Range("S25").Value = 562
$ws.Range("T25").Value = 16

# Row 26
$ws.Range("D26").Value = 44335
$ws.Range("K26").Value = 'Angeleno'
$ws.Range("M26").Value = 100
$ws.Range("Q26").Value = '$/bandeja 18 kilos granel'
$ws.Range("S26").Value = 583
$ws.Range("T26").Value = 18

# Row 27
$ws.Range("D27").Value = 44335
$ws.Range("K27").Value = 'Angeleno'
$ws.Range("M27").Value = 50
$ws.Range("Q27").Value = '$/bandeja 18 kilos granel'
$ws.Range("S27").Value = 500
$ws.Range("T27").Value = 18

# Row 28
$ws.Range("D28").Value = 44215
$ws.Range("K28").Value = 'Black Amber'
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 11000
$ws.Range("P28").Value = 10500
$ws.Range("S28").Value = 656

# Row 29
$ws.Range("D29").Value = 44215
$ws.Range("K29").Value = 'Black Amber'

# Row 30
$ws.Range("D30").Value = 44257
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 11000
$ws.Range("P30").Value = 10500
$ws.Range("Q30").Value = '$/caja 15 kilos granel'
$ws.Range("S30").Value = 700
$ws.Range("T30").Value = 15

# Row 31
$ws.Range("D31").Value = 44257
$ws.Range("N31").Value = 9000
$ws.Range("O31").Value = 9000
$ws.Range("P31").Value = 9000
$ws.Range("Q31").Value = '$/caja 15 kilos granel'
$ws.Range("S31").Value = 600
$ws.Range("T31").Value = 15

# Row 32
$ws.Range("D32").Value = 44251
$ws.Range("Q32").Value = '$/caja 16 kilos granel'
$ws.Range("S32").Value = 594
$ws.Range("T32").Value = 16

# Row 33
$ws.Range("D33").Value = 44251
$ws.Range("Q33").Value = '$/caja 16 kilos granel'
$ws.Range("S33").Value = 500
$ws.Range("T33").Value = 16

# Row 34
$ws.Range("D34").Value = 44279
$ws.Range("K34").Value = 'Black Amber'
$ws.Range("N34").Value = 9000
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 9500
$ws.Range("S34").Value = 528

# Row 35
$ws.Range("D35").Value = 44279
$ws.Range("K35").Value = 'Black Amber'
$ws.Range("N35").Value = 8000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 8000
$ws.Range("S35").Value = 444
